$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: update market price / profit figures
$ws.Range("H15").Value = 3047.3708
$ws.Range("I15").Value = 3047.3708
$ws.Range("K15").Value = 9142.1124
$ws.Range("M15").Value = -8973.1124
# Row 70: update market price / profit figures
$ws.Range("H70").Value = 1621.96
$ws.Range("I70").Value = 1128.1666
$ws.Range("J70").Value = 1777.8948
$ws.Range("K70").Value = 3384.4998
$ws.Range("L70").Value = 5333.6844
$ws.Range("M70").Value = -3114.4998
$ws.Range("N70").Value = -5873.6844
# Row 73: update market price / profit figures
$ws.Range("H73").Value = 1621.96
$ws.Range("I73").Value = 1128.1666
$ws.Range("J73").Value = 1777.8948
$ws.Range("K73").Value = 3384.4998
$ws.Range("L73").Value = 5333.6844
$ws.Range("M73").Value = -2448.4998
$ws.Range("N73").Value = -7205.6844
# Row 80: update market price / profit figures
$ws.Range("H80").Value = 9704.25
$ws.Range("I80").Value = 7398.4
$ws.Range("J80").Value = 10472.866
$ws.Range("K80").Value = 22195.2
$ws.Range("L80").Value = 31418.598
$ws.Range("M80").Value = -21197.2
$ws.Range("N80").Value = -33414.598
# Row 83: update market price / profit figures
$ws.Range("H83").Value = 9704.25
$ws.Range("I83").Value = 7398.4
$ws.Range("J83").Value = 10472.866
$ws.Range("K83").Value = 66585.59999999999
$ws.Range("L83").Value = 94255.79399999999
$ws.Range("M83").Value = -61593.59999999999
$ws.Range("N83").Value = -104239.794
# Row 88: update market price / profit figures
$ws.Range("H88").Value = 4776.077
$ws.Range("I88").Value = 4216.2
$ws.Range("J88").Value = 5126
$ws.Range("K88").Value = 4216.2
$ws.Range("L88").Value = 5126
$ws.Range("M88").Value = -3810.2
$ws.Range("N88").Value = -5938
# Row 91: update market price / profit figures
$ws.Range("H91").Value = 4776.077
$ws.Range("I91").Value = 4216.2
$ws.Range("J91").Value = 5126
$ws.Range("K91").Value = 4216.2
$ws.Range("L91").Value = 5126
$ws.Range("M91").Value = -2812.2
$ws.Range("N91").Value = -7934
# Row 94: update market price / profit figures
$ws.Range("H94").Value = 5974.1665
$ws.Range("I94").Value = 5608.1816
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 5608.1816
$ws.Range("L94").Value = 10000
$ws.Range("M94").Value = -5157.1816
$ws.Range("N94").Value = -10902
# Row 98: update market price / profit figures
$ws.Range("H98").Value = 1471.3636
$ws.Range("I98").Value = 1471.3636
$ws.Range("K98").Value = 1471.3636
$ws.Range("M98").Value = 26.63640000000009
# Row 122: update market price / profit figures
$ws.Range("H122").Value = 1471.3636
$ws.Range("I122").Value = 1471.3636
$ws.Range("K122").Value = 4414.0908
$ws.Range("M122").Value = -1964.0908
# Row 131: update market price / profit figures
$ws.Range("H131").Value = 58825576
$ws.Range("I131").Value = 90910120
$ws.Range("J131").Value = 3916.6667
$ws.Range("K131").Value = 272730360
$ws.Range("L131").Value = 11750.0001
$ws.Range("M131").Value = -272725320
$ws.Range("N131").Value = -21830.0001
# Row 138: update market price / profit figures
$ws.Range("H138").Value = 2859921
$ws.Range("I138").Value = 1353.9354
$ws.Range("J138").Value = 5132115.5
$ws.Range("K138").Value = 4061.8062
$ws.Range("L138").Value = 15396346.5
$ws.Range("M138").Value = 1078.1938
$ws.Range("N138").Value = -15406626.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2: update market price / profit figures
$ws.Range("H2").Value = 2121.2778
$ws.Range("J2").Value = 1765.2
$ws.Range("L2").Value = 1765.2
$ws.Range("N2").Value = -1991.2
# Row 116: update market price / profit figures
$ws.Range("H116").Value = 2121.2778
$ws.Range("J116").Value = 1765.2
$ws.Range("L116").Value = 1765.2
$ws.Range("N116").Value = -6353.2

$ws = $wb.Worksheets.Item("BSM")
# Row 3: update market price / profit figures
$ws.Range("H3").Value = 2121.2778
$ws.Range("J3").Value = 1765.2
$ws.Range("L3").Value = 1765.2
$ws.Range("N3").Value = -1993.2
# Row 107: update market price / profit figures
$ws.Range("H107").Value = 1891.875
$ws.Range("I107").Value = 1703.1428
$ws.Range("J107").Value = 3213
$ws.Range("K107").Value = 1703.1428
$ws.Range("L107").Value = 3213
$ws.Range("M107").Value = 216.8571999999999
$ws.Range("N107").Value = -7053

$ws = $wb.Worksheets.Item("CUL")
# Row 80: update market price / profit figures
$ws.Range("H80").Value = 500
$ws.Range("I80").Value = 500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -564
# Row 83: update market price / profit figures
$ws.Range("H83").Value = 500
$ws.Range("I83").Value = 500
$ws.Range("K83").Value = 4500
$ws.Range("M83").Value = 180

$ws = $wb.Worksheets.Item("GSM")
# Row 97: update market price / profit figures
$ws.Range("H97").Value = 1481.6666
$ws.Range("I97").Value = 1330.9
$ws.Range("J97").Value = 1670.125
$ws.Range("K97").Value = 1330.9
$ws.Range("L97").Value = 1670.125
$ws.Range("M97").Value = -834.9000000000001
$ws.Range("N97").Value = -2662.125

$ws = $wb.Worksheets.Item("LTW")
# Row 68: update market price / profit figures
$ws.Range("H68").Value = 9477
$ws.Range("I68").Value = 12305.2
$ws.Range("J68").Value = 3820.6
$ws.Range("K68").Value = 12305.2
$ws.Range("L68").Value = 3820.6
$ws.Range("M68").Value = -11556.2
$ws.Range("N68").Value = -5318.6
# Row 71: update market price / profit figures
$ws.Range("H71").Value = 9477
$ws.Range("I71").Value = 12305.2
$ws.Range("J71").Value = 3820.6
$ws.Range("K71").Value = 61526
$ws.Range("L71").Value = 19103
$ws.Range("M71").Value = -57782
$ws.Range("N71").Value = -26591
# Row 82: update market price / profit figures
$ws.Range("H82").Value = 867.2727
$ws.Range("I82").Value = 778
$ws.Range("J82").Value = 941.6667
$ws.Range("K82").Value = 778
$ws.Range("L82").Value = 941.6667
$ws.Range("M82").Value = -417
$ws.Range("N82").Value = -1663.6667
# Row 85: update market price / profit figures
$ws.Range("H85").Value = 867.2727
$ws.Range("I85").Value = 778
$ws.Range("J85").Value = 941.6667
$ws.Range("K85").Value = 778
$ws.Range("L85").Value = 941.6667
$ws.Range("M85").Value = 470
$ws.Range("N85").Value = -3437.6667
# Row 136: update market price / profit figures
$ws.Range("H136").Value = 11906542
$ws.Range("I136").Value = 25642674
$ws.Range("J136").Value = 1894.6666
$ws.Range("K136").Value = 76928022
$ws.Range("L136").Value = 5683.9998
$ws.Range("M136").Value = -76925472
$ws.Range("N136").Value = -10783.9998

$ws = $wb.Worksheets.Item("WVR")
# Row 62: update market price / profit figures
$ws.Range("H62").Value = 2350
$ws.Range("I62").Value = 2066.6667
$ws.Range("K62").Value = 2066.6667
$ws.Range("M62").Value = -1442.6667
# Row 65: update market price / profit figures
$ws.Range("H65").Value = 2350
$ws.Range("I65").Value = 2066.6667
$ws.Range("K65").Value = 10333.3335
$ws.Range("M65").Value = -7213.333500000001
# Row 132: update market price / profit figures
$ws.Range("H132").Value = 3639713.2
$ws.Range("I132").Value = 5408998
$ws.Range("J132").Value = 2850.7778
$ws.Range("K132").Value = 16226994
$ws.Range("L132").Value = 8552.3334
$ws.Range("M132").Value = -16224464
$ws.Range("N132").Value = -13612.3334
# Row 136: update market price / profit figures
$ws.Range("H136").Value = 50003764
$ws.Range("I136").Value = 66671252
$ws.Range("J136").Value = 1299.8
$ws.Range("K136").Value = 200013756
$ws.Range("L136").Value = 3899.4
$ws.Range("M136").Value = -200011206
$ws.Range("N136").Value = -8999.4
